$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated column F ("dSF") values re-pulled from the source data.
$updates = @{
    2  = -1
    4  = -3
    5  = 4
    6  = 2
    7  = 5
    8  = -3
    9  = -1
    10 = -2
    11 = 2
    12 = -7
    13 = -1
    14 = 1
    15 = 5
    16 = -4
    17 = 1
    18 = 0
    19 = 2
    20 = 4
    23 = -3
    24 = -1
    26 = 1
    29 = -2
    32 = 1
    33 = 1
    35 = 8
    36 = 2
    37 = -1
    38 = -9
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
